$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.032.79"
$ws.Range("E2").Value = "  +0.14%  "

$ws.Range("D3").Value = "1.844.76"
$ws.Range("E3").Value = "  -0.89%  "

$ws.Range("D4").Value = "'1.006"
$ws.Range("E4").Value = "  +0.43%  "

$ws.Range("D5").Value = "'329.94"
$ws.Range("E5").Value = "  -1.70%  "

$ws.Range("D6").Value = "'1.009"
$ws.Range("E6").Value = "  +0.81%  "

$ws.Range("D7").Value = "'0.4522"
$ws.Range("E7").Value = "  -3.73%  "

$ws.Range("D8").Value = "'0.3878"
$ws.Range("E8").Value = "  -0.61%  "

$ws.Range("D9").Value = "'47.48"
$ws.Range("E9").Value = "  +1.57%  "

$ws.Range("D10").Value = "'0.07724"
$ws.Range("E10").Value = "  -2.94%  "

$ws.Range("D11").Value = "'0.9696"
$ws.Range("E11").Value = "  -1.21%  "

$ws.Range("D12").Value = "'21.16"
$ws.Range("E12").Value = "  -1.59%  "

$ws.Range("D13").Value = "1.844.76"
$ws.Range("E13").Value = "  -0.72%  "

$ws.Range("D14").Value = "'5.754"
$ws.Range("E14").Value = "  -3.06%  "

$ws.Range("D15").Value = "'6.931"
$ws.Range("E15").Value = "  -3.62%  "

$ws.Range("D16").Value = "'1.009"
$ws.Range("E16").Value = "  +0.53%  "

$ws.Range("D17").Value = "'86.87"
$ws.Range("E17").Value = "  -4.97%  "

$ws.Range("D18").Value = "'0.06527"
$ws.Range("E18").Value = "  -1.33%  "

$ws.Range("D19").Value = "'0.00001014"
$ws.Range("E19").Value = "  -2.46%  "

$ws.Range("D20").Value = "'16.86"
$ws.Range("E20").Value = "  -3.53%  "

$ws.Range("D21").Value = "'1.016"
$ws.Range("E21").Value = "  +1.39%  "

$ws.Range("D22").Value = "27.992.32"
$ws.Range("E22").Value = "  -0.07%  "

$ws.Range("D23").Value = "'5.260"
$ws.Range("E23").Value = "  -2.48%  "

$ws.Range("D24").Value = "'10.52"
$ws.Range("E24").Value = "  -3.76%  "

$ws.Range("D25").Value = "'2.253"
$ws.Range("E25").Value = "  -1.43%  "

$ws.Range("D26").Value = "2.065.45"
$ws.Range("E26").Value = "  +0.00%  "

$ws.Range("D27").Value = "'156.19"
$ws.Range("E27").Value = "  -1.89%  "

$ws.Range("D28").Value = "'19.07"
$ws.Range("E28").Value = "  -2.35%  "

$ws.Range("D29").Value = "'2.017"
$ws.Range("E29").Value = "  -3.97%  "

$ws.Range("D30").Value = "'5.217"
$ws.Range("E30").Value = "  -4.43%  "

$ws.Range("D31").Value = "'116.22"
$ws.Range("E31").Value = "  -2.54%  "

$ws.Range("D32").Value = "'0.09228"
$ws.Range("E32").Value = "  -2.58%  "

$ws.Range("D33").Value = "'0.9274"
$ws.Range("E33").Value = "  -3.40%  "

$ws.Range("D34").Value = "'3.620"
$ws.Range("E34").Value = "  +1.20%  "

$ws.Range("E35").Value = "  +1.22%  "

$ws.Range("D36").Value = "'5.150"
$ws.Range("E36").Value = "  -2.84%  "

$ws.Range("D37").Value = "'0.05983"
$ws.Range("E37").Value = "  -1.58%  "

$ws.Range("D38").Value = "'0.02174"
$ws.Range("E38").Value = "  -3.50%  "

$ws.Range("D39").Value = "'8.062"
$ws.Range("E39").Value = "  -2.81%  "

$ws.Range("D40").Value = "'1.158"
$ws.Range("E40").Value = "  -0.36%  "

$ws.Range("D41").Value = "'1.006"
$ws.Range("E41").Value = "  +0.46%  "

$ws.Range("D42").Value = "'0.5606"
$ws.Range("E42").Value = "  -5.34%  "

$ws.Range("D43").Value = "'0.1776"
$ws.Range("E43").Value = "  -4.80%  "

$ws.Range("D44").Value = "'9.826"
$ws.Range("E44").Value = "  -3.65%  "

$ws.Range("D45").Value = "'1.241"
$ws.Range("E45").Value = "  -2.97%  "

$ws.Range("D46").Value = "'2.256"
$ws.Range("E46").Value = "  +22.80%  "

$ws.Range("D47").Value = "'0.07184"
$ws.Range("E47").Value = "  +4.48%  "

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'11.70"

$ws.Range("B49").Value = "Decentraland"
$ws.Range("C49").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D49").Value = "'0.5318"
$ws.Range("E49").Value = "  -4.08%  "

$ws.Range("D50").Value = "'1.861"
$ws.Range("E50").Value = "  -4.56%  "

$ws.Range("D51").Value = "'108.79"
$ws.Range("E51").Value = "  -2.39%  "
